# Update "想去人数" (want-to-go count) values in column F for the sheets
# that hold the 漫展 (convention) listing data: "展览" and "全部类型".
# Both sheets share identical data in this workbook, so the same set of
# row/value updates is applied to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 11676
    3  = 11285
    6  = 1021
    9  = 44
    11 = 10734
    12 = 4149
    14 = 6
    15 = 11
    17 = 1051
    21 = 11135
    22 = 10911
    24 = 28
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
